$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165025115013123
$ws.Range("B1").Value = 2.210657119750977
$ws.Range("C1").Value = 4.506183624267578
$ws.Range("D1").Value = 2.673350811004639
$ws.Range("E1").Value = 1.233798027038574
